## Negative Selection Prediction design doc - content updates
## Strategy: use Find.Execute to locate the legacy text, then replace that
## exact range (re-derived as a fresh Range so InsertXML actually overwrites
## it instead of inserting alongside it) with freshly built <w:r> runs via
## Range.InsertXML. This lets us reproduce the multi-run splits (and the
## <w:proofErr/> spell-check bookmarks) the author's copy of Word produced,
## while leaving the enclosing <w:p>/<w:pPr> (and its w14:paraId etc.)
## completely untouched.

$d = $word.ActiveDocument

function Escape-Xml([string]$s) {
    $s = $s -replace "&", "&amp;"
    $s = $s -replace "<", "&lt;"
    $s = $s -replace ">", "&gt;"
    return $s
}

function Build-Run([string]$text) {
    $preserve = ""
    if ($text.Length -eq 0 -or $text[0] -eq " " -or $text[$text.Length - 1] -eq " ") {
        $preserve = ' xml:space="preserve"'
    }
    $escaped = Escape-Xml $text
    return "<w:r><w:rPr><w:rFonts w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`" w:cs=`"Times New Roman`"/><w:sz w:val=`"22`"/><w:szCs w:val=`"22`"/></w:rPr><w:t$preserve>$escaped</w:t></w:r>"
}

# Each $segments entry is either a plain string (-> one run) or a hashtable
# @{ Text = "..."; SpellErr = $true } to additionally wrap the run with
# <w:proofErr w:type="spellStart"/> ... <w:proofErr w:type="spellEnd"/>,
# matching what Word's background spell checker stamped into the XML.
function Build-Runs([object[]]$segments) {
    $body = ""
    foreach ($seg in $segments) {
        if ($seg -is [hashtable]) {
            $run = Build-Run $seg.Text
            if ($seg.ContainsKey("SpellErr") -and $seg.SpellErr) {
                $body += '<w:proofErr w:type="spellStart"/>' + $run + '<w:proofErr w:type="spellEnd"/>'
            } else {
                $body += $run
            }
        } else {
            $body += Build-Run $seg
        }
    }
    return $body
}

function Build-PkgXml([string]$pBody) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $pBody + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

function Replace-FoundTextWithRuns([string]$findText, [object[]]$segments) {
    $hit = $d.Content
    $found = $hit.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find text: $findText"
    }
    $target = $d.Range($hit.Start, $hit.End)
    $xml = Build-PkgXml (Build-Runs $segments)
    $target.InsertXML($xml)
}

# 1) "Input: ..." paragraph -> split into 4 runs describing the new FASTA input.
Replace-FoundTextWithRuns `
    "Input: Data about various undifferentiated T Cells in the proper environment conditions to have them undergo selection." `
    @(
        "Input: ",
        "FASTA sequences of peptide antigen, MHC I and II binding region, and T",
        "-",
        "Cell antigen binding region."
    )

# 2) "Output: ..." paragraph -> split into many runs describing survival based
#    on free-energy comparison; "TCell" gets wrapped in proofErr markers.
Replace-FoundTextWithRuns `
    "Output: List of T-Cells that survived negative selection and which ones survived and which ones didn’t. The chance of survival is based on a p value that results after determining the strength of binding to AIRE and then comparing with the average strength of binding of other T-Cells. " `
    @(
        "Output: ",
        "Whether a specific",
        " T-Cell survived negative selection",
        ".",
        " ",
        "S",
        "urvival is based on ",
        "change in free energy after the ",
        @{ Text = "TCell"; SpellErr = $true },
        " binds to the peptide antigen and the MHC is bound to the peptide. In comparing these free energies we can conclude that the lowest one has stronger affinity for the peptide",
        "."
    )

# 3) Remove the "Raw Data of undifferentiated T-Cells..." paragraph entirely
#    (the bullet content right after "Requirements:").
$hit = $d.Content
$found = $hit.Find.Execute("Raw Data of undifferentiated T-Cells to determine their potential binding coefficient. Afterwards setting a limit to how strongly it can be bound to AIRE.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the Raw Data paragraph"
}
$para = $d.Range($hit.Start, $hit.End).Paragraphs(1)
$para.Range.Delete()

# 4) The (now) empty paragraph that used to trail the deleted bullet gains new
#    content describing the Fasta file data source; "Fasta" is wrapped with
#    proofErr markers as Word's spell checker flagged it.
$hit = $d.Content
$found = $hit.Find.Execute("Great Lakes Research cluster.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    throw "Fasta sentence already present before edit"
}

# Locate the trailing empty paragraph right before the section break (last
# paragraph in the body) and fill it in.
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$target = $lastPara.Range
$target = $d.Range($target.Start, $target.End - 1)
$xml = Build-PkgXml (Build-Runs @(
        @{ Text = "Fasta"; SpellErr = $true },
        " file containing the AA sequences. Great Lakes Research cluster."
    ))
$target.InsertXML($xml)
